$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting rows 17-225 down to 18-226.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with a fresh data record (same
# market/category/variety/quality as the rest of the sheet, with its own
# date and volume).
$ws.Cells.Item(17, 1).Value = 3
$ws.Cells.Item(17, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44545
$ws.Cells.Item(17, 5).Value = 5
$ws.Cells.Item(17, 6).Value = 100112039
$ws.Cells.Item(17, 7).Value = "Ciboulette"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 180
$ws.Cells.Item(17, 11).Value = 1500
$ws.Cells.Item(17, 12).Value = 1500
$ws.Cells.Item(17, 13).Value = 1500
$ws.Cells.Item(17, 14).Value = "$/docena de atados"
$ws.Cells.Item(17, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 16).Value = 500
$ws.Cells.Item(17, 17).Value = 3
$ws.Cells.Item(17, 18).Value = "Hortaliza"
